# Add a new "colour" column (D) to the Brews sheet.
# The header goes in D1, and rows that don't already have a colour value
# get a blank cell styled with a flat fill colour instead of text (this is
# how the importer "parses" a colour out of a blank, styled data-row cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the colour column.
$ws.Range("D1").Value = "colour"

# Row 2 (hop hog / ipa) -> orange fill (#FFC000), cell itself stays blank.
# Interior.Color uses the COM BGR-packed long: R + G*256 + B*65536.
$ws.Range("D2").Interior.Color = 49407    # 0x00, 0xC0, 0xFF -> R=255 G=192 B=0

# Row 3 (fanta pants / american amber) -> pink/magenta fill (#C600AE), blank cell.
$ws.Range("D3").Interior.Color = 11403462 # R=198 G=0 B=174

# Leave the active selection on the last edited cell, matching the saved view state.
$ws.Range("D3").Select() | Out-Null
